# The workbook tracks daily produce-market prices for "Choclo" (corn) at
# Terminal La Palmera de La Serena. A new daily record was added for the
# market on 2023-08-04 (serial date 45142). It was inserted as a new row
# at row 777, pushing the existing rows 777-889 down to 778-890, and it
# reuses the exact same data values the (then) row 777 had - only the
# date changed.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the current row 777, shifting rows 777:889 down
# to 778:890.
$ws.Range("A777:R777").EntireRow.Insert()

# Populate the newly inserted row 777 with the new record.
$ws.Range("A777").Value = 8
$ws.Range("B777").Value = "Terminal La Palmera de La Serena"
$ws.Range("C777").Value = "Coquimbo"
$ws.Range("D777").Value = 45142
$ws.Range("E777").Value = 4
$ws.Range("F777").Value = 100112024
$ws.Range("G777").Value = "Choclo"
$ws.Range("H777").Value = "Dulce o Americano"
$ws.Range("I777").Value = "Primera"
$ws.Range("J777").Value = 580
$ws.Range("K777").Value = 37000
$ws.Range("L777").Value = 38000
$ws.Range("M777").Value = 37500
$ws.Range("N777").Value = "`$/malla 70 unidades"
$ws.Range("O777").Value = "Región de Arica y Parinacota"
$ws.Range("P777").Value = 536
$ws.Range("Q777").Value = 70
$ws.Range("R777").Value = "Hortaliza"
